# Insert two new weekly-report rows at the top of the Camote price table
# (row 74), pushing the existing data (old rows 74-98) down to rows 76-100.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("74:75").Insert()

# New row 74
$ws.Range("A74").Value = 9
$ws.Range("B74").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C74").Value = 'Metropolitana'
$ws.Range("D74").Value = 44795
$ws.Range("E74").Value = 13
$ws.Range("F74").Value = 100114002
$ws.Range("G74").Value = 'Camote'
$ws.Range("H74").Value = 'Sin especificar'
$ws.Range("I74").Value = 'Primera'
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 13000
$ws.Range("L74").Value = 14000
$ws.Range("M74").Value = 13471
$ws.Range("N74").Value = '$/caja 18 kilos'
$ws.Range("O74").Value = 'Perú'
$ws.Range("P74").Value = 748
$ws.Range("Q74").Value = 18
$ws.Range("R74").Value = 'Hortaliza'

# New row 75
$ws.Range("A75").Value = 9
$ws.Range("B75").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C75").Value = 'Metropolitana'
$ws.Range("D75").Value = 44795
$ws.Range("E75").Value = 13
$ws.Range("F75").Value = 100114002
$ws.Range("G75").Value = 'Camote'
$ws.Range("H75").Value = 'Sin especificar'
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 1400
$ws.Range("K75").Value = 13000
$ws.Range("L75").Value = 14000
$ws.Range("M75").Value = 13571
$ws.Range("N75").Value = '$/malla 18 kilos'
$ws.Range("O75").Value = 'Perú'
$ws.Range("P75").Value = 754
$ws.Range("Q75").Value = 18
$ws.Range("R75").Value = 'Hortaliza'
